$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    0.22018249879080543,
    -0.0059999999653719271,
    0.017144659581004973,
    -0.0079999999443476355,
    -0.0029999999685212941,
    -0.0019999999661113321,
    -0.0099999999221798141,
    -0.0099999999217494917,
    0.023596668119048481,
    -0.0019999999647950517,
    -0.0029999999590071269,
    -0.0034999999557396855,
    -0.0034999999541760474,
    -0.0079999999290665258,
    -0.00099999996776212186,
    0.0031682341005225823,
    -0.0019999999617965614,
    -0.0039999999505528905,
    -0.0039999999755546689,
    -0.051966594094952967,
    -0.0039999999663242747,
    -0.0039999999661128882,
    -0.0049999999624787961,
    -0.019999999876278096,
    -0.01999999987465273,
    -0.0024999999620476387,
    -0.0024999999600501255,
    -0.0019999999561131077,
    -0.006999999922605582,
    -0.059999999626469513,
    0.0072717299436089178,
    -0.0099999999020621289,
    -0.003999999935198062
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
